$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.174200534820557
$ws.Range("B1").Value = 2.329868316650391
$ws.Range("C1").Value = 4.490339756011963
$ws.Range("D1").Value = 3.519518613815308
$ws.Range("E1").Value = 1.216854572296143
